$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits right after the
#    "cd ~/servers/wso2esb-4.8.1" command text. It needs to move
#    down to the end of the "Please take a reasonably good look at
#    this." bullet (after the sentence that gets added below).
#    Remove it from its old spot first.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Append the reminder sentence right after "Please take a
#    reasonably good look at this."
# ------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("Please take a reasonably good look at this.", $false, $false, $false, $false, $false, $true, 1, $false, "Please take a reasonably good look at this. Also, make sure your Application Server and ESB are running.", 2)

# ------------------------------------------------------------------
# 3. Re-create the "_GoBack" bookmark right after the sentence we
#    just inserted (and before the trailing line break that ends
#    the bullet).
# ------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("Also, make sure your Application Server and ESB are running.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $find2.Parent
$target.Collapse(0)
$d.Bookmarks.Add("_GoBack", $target)
